$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Simple single-value cell updates
$t.Cell(1,1).Range.Text = "0M"
$t.Cell(2,1).Range.Text = "0M"
$t.Cell(3,1).Range.Text = "0M"
$t.Cell(4,1).Range.Text = "111"
$t.Cell(6,1).Range.Text = "0.00048"
$t.Cell(7,1).Range.Text = "0.00017"
$t.Cell(9,1).Range.Text = "0.00031"
$t.Cell(10,1).Range.Text = "0.00035"
$t.Cell(11,1).Range.Text = "0.00042"
$t.Cell(12,1).Range.Text = "0.02197"

# Collapse the multi-run, tab-separated cells down to a single value
$t.Cell(44,1).Range.Text = "99.96"
$t.Cell(45,1).Range.Text = "0.02"
$t.Cell(46,1).Range.Text = "57"
